# 自动更新Excel文件 - 2026-01-17 23:13:10
# Decrement "剩余" (E column) by 1 for each row, representing one day elapsed.
# When 剩余 reaches 1, it wraps back to 10 and 开始时间 (F column, a new cycle
# start date stored as YYYYMMDD) advances by 10 days. Row 36's date value is
# malformed in the source data and is left untouched, matching the original edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 20260112
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 20260112
$ws.Range("E4").Value = 8
$ws.Range("F4").Value = 20260112
$ws.Range("E5").Value = 6
$ws.Range("F5").Value = 20260114
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 20260112
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 20260114
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 20260112
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = 20260114
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 20260112
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = 20260112
$ws.Range("E12").Value = 6
$ws.Range("F12").Value = 20260114
$ws.Range("E13").Value = 8
$ws.Range("F13").Value = 20260112
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 20260112
$ws.Range("E15").Value = 8
$ws.Range("F15").Value = 20260112
$ws.Range("E16").Value = 10
$ws.Range("F16").Value = 20260118
$ws.Range("E17").Value = 6
$ws.Range("F17").Value = 20260114
$ws.Range("E18").Value = 9
$ws.Range("F18").Value = 20260117
$ws.Range("E19").Value = 9
$ws.Range("F19").Value = 20260117
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = 20260117
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = 20260117
$ws.Range("E22").Value = 6
$ws.Range("F22").Value = 20260114
$ws.Range("E23").Value = 6
$ws.Range("F23").Value = 20260114
$ws.Range("E24").Value = 6
$ws.Range("F24").Value = 20260114
$ws.Range("E25").Value = 6
$ws.Range("F25").Value = 20260114
$ws.Range("E26").Value = 6
$ws.Range("F26").Value = 20260114
$ws.Range("E27").Value = 2
$ws.Range("F27").Value = 20260113
$ws.Range("E28").Value = 9
$ws.Range("F28").Value = 20260117
$ws.Range("E29").Value = 9
$ws.Range("F29").Value = 20260117
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 20260117
$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 20260117
$ws.Range("E32").Value = 9
$ws.Range("F32").Value = 20260117
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 20260117
$ws.Range("E34").Value = 9
$ws.Range("F34").Value = 20260117
$ws.Range("E35").Value = 9
$ws.Range("F35").Value = 20260117
$ws.Range("E37").Value = 9
$ws.Range("F37").Value = 20260117
$ws.Range("E38").Value = 9
$ws.Range("F38").Value = 20260117
$ws.Range("E39").Value = 9
$ws.Range("F39").Value = 20260117
$ws.Range("E40").Value = 1
$ws.Range("F40").Value = 20260112
$ws.Range("E41").Value = 1
$ws.Range("F41").Value = 20260112
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 20260117
$ws.Range("E43").Value = 6
$ws.Range("F43").Value = 20260114
$ws.Range("E44").Value = 1
$ws.Range("F44").Value = 20260112
$ws.Range("E45").Value = 6
$ws.Range("F45").Value = 20260114
$ws.Range("E46").Value = 1
$ws.Range("F46").Value = 20260112
$ws.Range("E47").Value = 9
$ws.Range("F47").Value = 20260117
$ws.Range("E48").Value = 1
$ws.Range("F48").Value = 20260112
$ws.Range("E49").Value = 2
$ws.Range("F49").Value = 20260113
$ws.Range("E50").Value = 4
$ws.Range("F50").Value = 20260112
$ws.Range("E51").Value = 4
$ws.Range("F51").Value = 20260112
$ws.Range("E52").Value = 4
$ws.Range("F52").Value = 20260112
$ws.Range("E53").Value = 4
$ws.Range("F53").Value = 20260112
$ws.Range("E54").Value = 4
$ws.Range("F54").Value = 20260112
$ws.Range("E55").Value = 4
$ws.Range("F55").Value = 20260112
$ws.Range("E56").Value = 4
$ws.Range("F56").Value = 20260112
$ws.Range("E57").Value = 4
$ws.Range("F57").Value = 20260112
$ws.Range("E58").Value = 8
$ws.Range("F58").Value = 20260116
$ws.Range("E59").Value = 8
$ws.Range("F59").Value = 20260116
$ws.Range("E60").Value = 8
$ws.Range("F60").Value = 20260116
$ws.Range("E61").Value = 2
$ws.Range("F61").Value = 20260113
$ws.Range("E62").Value = 8
$ws.Range("F62").Value = 20260116
$ws.Range("E63").Value = 8
$ws.Range("F63").Value = 20260116
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = 20260116
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 20260117
$ws.Range("E66").Value = 9
$ws.Range("F66").Value = 20260117
$ws.Range("E67").Value = 9
$ws.Range("F67").Value = 20260117
$ws.Range("E68").Value = 9
$ws.Range("F68").Value = 20260117
$ws.Range("E69").Value = 9
$ws.Range("F69").Value = 20260117
$ws.Range("E70").Value = 10
$ws.Range("F70").Value = 20260118
$ws.Range("E71").Value = 10
$ws.Range("F71").Value = 20260118
$ws.Range("E72").Value = 10
$ws.Range("F72").Value = 20260118
$ws.Range("E73").Value = 10
$ws.Range("F73").Value = 20260118
$ws.Range("E74").Value = 10
$ws.Range("F74").Value = 20260118
$ws.Range("E75").Value = 10
$ws.Range("F75").Value = 20260118
$ws.Range("E76").Value = 10
$ws.Range("F76").Value = 20260118
$ws.Range("E77").Value = 3
$ws.Range("F77").Value = 20260111
$ws.Range("E78").Value = 3
$ws.Range("F78").Value = 20260111
$ws.Range("E79").Value = 3
$ws.Range("F79").Value = 20260111
$ws.Range("E80").Value = 3
$ws.Range("F80").Value = 20260111
$ws.Range("E81").Value = 3
$ws.Range("F81").Value = 20260111
$ws.Range("E82").Value = 3
$ws.Range("F82").Value = 20260111
$ws.Range("E83").Value = 3
$ws.Range("F83").Value = 20260111
$ws.Range("E84").Value = 3
$ws.Range("F84").Value = 20260111
$ws.Range("E85").Value = 3
$ws.Range("F85").Value = 20260111
$ws.Range("E86").Value = 3
$ws.Range("F86").Value = 20260111
$ws.Range("E87").Value = 1
$ws.Range("F87").Value = 20260112
$ws.Range("E88").Value = 1
$ws.Range("F88").Value = 20260112
$ws.Range("E89").Value = 1
$ws.Range("F89").Value = 20260112
$ws.Range("E90").Value = 1
$ws.Range("F90").Value = 20260112
$ws.Range("E91").Value = 6
$ws.Range("F91").Value = 20260114
$ws.Range("E92").Value = 1
$ws.Range("F92").Value = 20260112
$ws.Range("E93").Value = 3
$ws.Range("F93").Value = 20260111
$ws.Range("E94").Value = 4
$ws.Range("F94").Value = 20260115
$ws.Range("E95").Value = 2
$ws.Range("F95").Value = 20260110
$ws.Range("E96").Value = 10
$ws.Range("F96").Value = 20260118
$ws.Range("E97").Value = 10
$ws.Range("F97").Value = 20260118
$ws.Range("E98").Value = 10
$ws.Range("F98").Value = 20260118
$ws.Range("E99").Value = 10
$ws.Range("F99").Value = 20260118
